# Add 2022-Q3 data:
#  1. Insert a new "2022-Q3" worksheet right before "2022-Q2" (so the tab
#     order becomes: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3,
#     2021-Q2, 2021-Q1).
#  2. Insert a new row into "总计" for the 2022-Q3 summary figures and push
#     the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned before "2022-Q2".
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$q3Rows = @(
    @(0, "159792", "富国中证港股通互联网ETF", "18.08", "99.26", "2.99", "0.5406", 10),
    @(1, "513770", "华宝中证港股通互联网ETF", "3.80", "98.21", "2.94", "0.1117", 10)
)

foreach ($row in $q3Rows) {
    $r = [int]$row[0] + 2
    $q3Sheet.Cells.Item($r, 1).Value = $row[0]
    $q3Sheet.Cells.Item($r, 2).Value = $row[1]
    $q3Sheet.Cells.Item($r, 3).Value = $row[2]
    $q3Sheet.Cells.Item($r, 4).Value = $row[3]
    $q3Sheet.Cells.Item($r, 5).Value = $row[4]
    $q3Sheet.Cells.Item($r, 6).Value = $row[5]
    $q3Sheet.Cells.Item($r, 7).Value = $row[6]
    $q3Sheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add the 2022-Q3 row and push the
#    other quarters down by one row (2022-Q2 -> row 3, ..., 2021-Q1 -> row 8).
#    Writing every row's final value directly (rather than using Insert(),
#    which drags along stray blank-cell formatting) keeps the styling
#    identical to the source rows.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(2, 0, "2022-Q3", 2, 0.65),
    @(3, 1, "2022-Q2", 5, 0.78),
    @(4, 2, "2022-Q1", 1, 0.07000000000000001),
    @(5, 3, "2021-Q4", 6, 1.6),
    @(6, 4, "2021-Q3", 10, 4.19),
    @(7, 5, "2021-Q2", 10, 4.22),
    @(8, 6, "2021-Q1", 7, 2.95)
)

foreach ($row in $summaryRows) {
    $r = $row[0]
    $summary.Cells.Item($r, 1).Value = $row[1]
    $summary.Cells.Item($r, 2).Value = $row[2]
    $summary.Cells.Item($r, 3).Value = $row[3]
    $summary.Cells.Item($r, 4).Value = $row[4]
}

# Row 8 sits beyond the original A1:D7 range, so its index cell (A8) starts
# out unstyled; copy the formatting already used by the other index cells
# in column A (e.g. A7) so it matches exactly.
$summary.Range("A7").Copy($summary.Range("A8"))
$summary.Cells.Item(8, 1).Value = 6
